$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value  = -8.669
$ws.Range("D10").Value = -8.099
$ws.Range("D12").Value = -6.725
$ws.Range("D18").Value = -8.318999999999999
$ws.Range("D37").Value = -8.315
$ws.Range("D55").Value = -8.218
$ws.Range("D68").Value = -7.229000000000001
$ws.Range("D77").Value = -7.840000000000001
$ws.Range("D78").Value = -8.279
$ws.Range("D81").Value = -7.754
$ws.Range("D82").Value = -8.17
